# Apply updated KNN imputation result values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D7"  = -7.520999999999999
    "C9"  = -10.857
    "D12" = -7.375999999999999
    "E15" = 16.131
    "C18" = -12.054
    "C20" = -12.205
    "D26" = -8.062000000000001
    "C27" = -13.111
    "D27" = -8.452999999999999
    "D29" = -7.292
    "D37" = -7.943000000000001
    "D38" = -7.593999999999999
    "E38" = 16.692
    "E44" = 16.735
    "D51" = -8.175000000000001
    "E51" = 16.724
    "D55" = -7.918000000000001
    "E57" = 16.616
    "E63" = 17.731
    "C69" = -11.199
    "D69" = -7.219999999999999
    "D70" = -6.9
    "E70" = 17.849
    "C76" = -13.032
    "C82" = -11.991
    "D83" = -8.373000000000001
    "E99" = 16.583
    "D102" = -7.861
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
